$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.623.37"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.397.02"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "562.98"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "140.93"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").Value = "2.403.30"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "5.15"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "60.275.65"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "2.402.76"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "7.93"
$ws.Range("E19").Value = "  +6.23%  "
$ws.Range("D20").Value = "10.65"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "323.59"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "6.08"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("D26").Value = "64.98"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "563.03"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("E28").Value = "  -4.54%  "
$ws.Range("D29").Value = "2.513.02"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").Value = "8.09"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("E32").Value = "  -1.65%  "
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("E36").Value = "  +3.44%  "
$ws.Range("D37").Value = "152.22"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "5.13"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").Value = "2.49"
$ws.Range("E45").Value = "  +6.34%  "
$ws.Range("D46").Value = "0.0₆0284"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "141.33"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").Value = "19.23"
$ws.Range("E51").Value = "  -1.62%  "
